$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. testdata_Mean sheet (sheet1): row 28 gains Area1 / NA values,
#    selection moves to A29:XFD29.
# ------------------------------------------------------------------
$wsMean = $wb.Worksheets.Item("testdata_Mean")
$wsMean.Range("A28").Value = "Area1"
$wsMean.Range("B28").Value = "NA"
$wsMean.Range("A29:XFD29").Select()

# ------------------------------------------------------------------
# 2. Add the new results sheet "testdata_Mean_results_NA" after the
#    existing "testdata_Mean_results" sheet.
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "testdata_Mean_results_NA"

# Re-fetch stable handles by name right before the move - positional
# references captured before an Add()/Move() go stale once the sheet
# collection is reshuffled.
$wsNA = $wb.Worksheets.Item("testdata_Mean_results_NA")
$wsResults = $wb.Worksheets.Item("testdata_Mean_results")
$wsNA.Move($null, $wsResults)
$wsNA = $wb.Worksheets.Item("testdata_Mean_results_NA")

# Header row (bold style, same headings as testdata_Mean_results)
$wsNA.Range("A1").Value = "area"
$wsNA.Range("B1").Value = "value_sum"
$wsNA.Range("C1").Value = "value_count"
$wsNA.Range("D1").Value = "stdev"
$wsNA.Range("E1").Value = "Value"
$wsNA.Range("F1").Value = "lower_95_ci"
$wsNA.Range("G1").Value = "upper_95_ci"
$wsNA.Range("H1").Value = "lower_99_8_ci"
$wsNA.Range("I1").Value = "upper_99_8_ci"
$wsNA.Range("J1").Value = "Statistic"
$wsNA.Range("K1").Value = "Method"
$wsNA.Range("A1:K1").Font.Bold = $true
$wsNA.Range("A1:K1").Interior.Pattern = -4124
$wsNA.Range("A1:K1").Interior.ThemeColor = 1
$wsNA.Range("A1:K1").Interior.TintAndShade = -0.14999847407452621

# Row 2 - Area1 group is now fully NA
$wsNA.Range("A2").Value = "Area1"
$wsNA.Range("B2").Value = "NA"
$wsNA.Range("C2").Value = 8
$wsNA.Range("D2").Value = "NA"
$wsNA.Range("E2").Value = "NA"
$wsNA.Range("F2").Value = "NA"
$wsNA.Range("G2").Value = "NA"
$wsNA.Range("H2").Value = "NA"
$wsNA.Range("I2").Value = "NA"
$wsNA.Range("J2").Value = "Mean"
$wsNA.Range("K2").Value = "Student's t-distribution"

# Row 3 - Area2 group unchanged from testdata_Mean_results
$wsNA.Range("A3").Value = "Area2"
$wsNA.Range("B3").Value = 102221.33323999999
$wsNA.Range("C3").Value = 18
$wsNA.Range("D3").Value = 2117.8317161590671
$wsNA.Range("E3").Value = 5678.9629577777778
$wsNA.Range("F3").Value = 4625.7900224529722
$wsNA.Range("G3").Value = 6732.1358931025834
$wsNA.Range("H3").Value = 3859.0770997295967
$wsNA.Range("I3").Value = 7498.8488158259588
$wsNA.Range("J3").Value = "Mean"
$wsNA.Range("K3").Value = "Student's t-distribution"

# Row 4 - totals row, B4/C4 are formulas (SUM ignores the NA text cells)
$wsNA.Range("A4").Value = "No grouping"
$wsNA.Range("B4").Formula = "=SUM(B2:B3)"
$wsNA.Range("C4").Formula = "=SUM(C2:C3)"
$wsNA.Range("D4").Value = 3171.8018122194453
$wsNA.Range("E4").Value = 3948.1141538461534
$wsNA.Range("F4").Value = 2666.9956767458489
$wsNA.Range("G4").Value = 5229.2326309464579
$wsNA.Range("H4").Value = 1801.9535385474737
$wsNA.Range("I4").Value = 6094.274769144833
$wsNA.Range("J4").Value = "Mean"
$wsNA.Range("K4").Value = "Student's t-distribution"

$wsNA.Range("A2").Select()
